$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: clear D3 (was a number, becomes blank)
$ws.Range("D3").ClearContents()

# Row 4: C4 value updated
$ws.Range("C4").Value = 0

# Row 5: C5 value updated
$ws.Range("C5").Value = 544.9647926184913

# Row 7: rename "Other" -> "Biogas" and update D7 value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 57.90768029780688

# New row 8: "Other" with D8 value, matching style of row 7
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Other"
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = 900.5104274015649
